$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "아프리카 최빈국 말라위 청년에게 자립 기회를!"
$ws.Range("A3").Value = "“집에 있으면 뭐하나... 여기 와야 재밌지.”"
$ws.Range("A4").Value = "우리는 함께 먹을수록 단단해진다"
$ws.Range("A5").Value = "철수의 꿈을 함께 찾아주세요."
$ws.Range("A6").Value = "우리 집을 만들어주세요!"
$ws.Range("A7").Value = "6.25참전 유공자분들께 지팡이를 전달해 주세요"
$ws.Range("A8").Value = "따뜻한 한 끼 밥상, 나눔으로 배부른 보통 일상"
$ws.Range("A9").Value = "“우리는 경로당 아니면 갈 데가 없어.”"
$ws.Range("A10").Value = "올바른 장애 첫인상을 함께 만들어요"
